$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - row => new F value
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 8
$wsExpo.Range("F4").Value = 13266
$wsExpo.Range("F6").Value = 3
$wsExpo.Range("F8").Value = 114
$wsExpo.Range("F9").Value = 113
$wsExpo.Range("F10").Value = 63
$wsExpo.Range("F13").Value = 13231
$wsExpo.Range("F14").Value = 327
$wsExpo.Range("F16").Value = 8850
$wsExpo.Range("F17").Value = 7930
$wsExpo.Range("F24").Value = 1006
$wsExpo.Range("F26").Value = 17
$wsExpo.Range("F29").Value = 101
$wsExpo.Range("F30").Value = 357

# Sheet "全部类型" (all types) - row => new F value
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 8
$wsAll.Range("F5").Value = 13266
$wsAll.Range("F7").Value = 3
$wsAll.Range("F9").Value = 114
$wsAll.Range("F10").Value = 113
$wsAll.Range("F11").Value = 63
$wsAll.Range("F14").Value = 13231
$wsAll.Range("F15").Value = 327
$wsAll.Range("F17").Value = 8851
$wsAll.Range("F18").Value = 7930
$wsAll.Range("F25").Value = 1006
$wsAll.Range("F27").Value = 17
$wsAll.Range("F32").Value = 101
$wsAll.Range("F33").Value = 357
